# "edited typo in flasher_v1_0"
#
# The Config sheet's firmware table (Tabelle1) and the Constants sheet's
# firmware-selection table (Tabelle2) get two extra rows each for a
# renamed/rebuilt set of DFU packages:
#   - master_dfu_package.zip (Number 45, was the "flasher" row)
#   - client_dfu_package.zip (Numbers 5 and 11)
#   - server_dfu_package.zip (Numbers 14 and 15)
# The now-unused "Group ID" column values on the Config sheet are cleared.

$wb = $excel.ActiveWorkbook
$wsConfig = $wb.Worksheets.Item("Config")
$wsConstants = $wb.Worksheets.Item("Constants")

# --- Grow the two tables so there is room for the new rows -----------------
$tblConfig = $wsConfig.ListObjects.Item("Tabelle1")
$tblConfig.ListRows.Add() | Out-Null
$tblConfig.ListRows.Add() | Out-Null

$tblFirmwareList = $wsConstants.ListObjects.Item("Tabelle2")
$tblFirmwareList.ListRows.Add() | Out-Null
$tblFirmwareList.ListRows.Add() | Out-Null
$tblFirmwareList.ListRows.Add() | Out-Null

# --- Constants sheet: new firmware-selection entries (A6:A8) --------------
$wsConstants.Range("A6").Value = "master_dfu_package.zip"
$wsConstants.Range("A7").Value = "client_dfu_package.zip"
$wsConstants.Range("A8").Value = "server_dfu_package.zip"

# --- Config sheet: update existing rows 2-4, clear Group ID -----------------
$wsConfig.Range("A2").Value = 45
$wsConfig.Range("B2").Value = "master_dfu_package.zip"
$wsConfig.Range("D2").ClearContents() | Out-Null

$wsConfig.Range("A3").Value = 5
$wsConfig.Range("B3").Value = "client_dfu_package.zip"
$wsConfig.Range("D3").ClearContents() | Out-Null

$wsConfig.Range("A4").Value = 11
$wsConfig.Range("B4").Value = "client_dfu_package.zip"
$wsConfig.Range("D4").ClearContents() | Out-Null

# --- Config sheet: fill the two new rows (5-6) ------------------------------
$wsConfig.Range("A5").Value = 14
$wsConfig.Range("B5").Value = "server_dfu_package.zip"
$wsConfig.Range("C5").Formula = "=VLOOKUP(Tabelle1[[#This Row],[Number]],Tabelle3[],2,FALSE)"

$wsConfig.Range("A6").Value = 15
$wsConfig.Range("B6").Value = "server_dfu_package.zip"
$wsConfig.Range("C6").Formula = "=VLOOKUP(Tabelle1[[#This Row],[Number]],Tabelle3[],2,FALSE)"

# --- Restore the active sheet / selection --------------------------------
$wsConstants.Activate()
$wsConstants.Range("A13").Select() | Out-Null

$wsConfig.Activate()
$wsConfig.Range("K11").Select() | Out-Null
